$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 40 (pushes old rows 40-48 down to 41-49),
# to host the new "Scanned: Awaiting results" record.
$ws.Rows.Item(40).Insert()

# Row 39 - update description text (147 -> 185 days context change is NOT in row 39,
# but the "valid scan records" -> "valid scan record" wording tweak is in E39)
$ws.Range("E39").Value = "The number of unique participants in the LDCT table with at least one valid scan record (dated and outcomed as being performed) but do not appear in the NCRAS dataset with a lung cancer diagnosis."

# Row 40 (new) - "Scanned: Awaiting results"
$ws.Range("A40").Value = "cancer_outcome"
$ws.Range("B40").Value = "Scanned: Awaiting results"
$ws.Range("C40").Value = "LDCT AND NCRAS consolidated cancer outcomes"
$ws.Range("D40").Value = "The number of participants who had at least one low-dose CT scan and do not appear in the NCRAS dataset with a diagnosis of lung cancer and where their scan took place after the period covered by NCRAS data."
$ws.Range("E40").Value = "The number of unique participants in the LDCT table with at least one valid scan record (dated and outcomed as being performed) but do not appear in the NCRAS dataset with a lung cancer diagnosis and where the scan took place after 27th Feb 2023 (185 days before the end of the NCRAS cancer diagnosis data) meaning we are unsure of the result of the scan."
$ws.Range("A40:E40").Style = $ws.Range("A39:E39").Style
$ws.Rows.Item(40).RowHeight = 75

# Row 41 (was 40) - 147 -> 185 days
$ws.Range("D41").Value = "The number of participants with a low-dose CT scan or were assessed as high risk at LHC and a lung cancer diagnosed within 185 days of their TLHC contact."
$ws.Range("E41").Value = "The number of participants who had either a) a low-dose CT scan or b) a LHC at which they were assessed as being high risk and eligible for a scan, and also have a lung cancer diagnosed within 185 days of their TLHC contact  in the NCRAS dataset."

# Row 43 (was 42) - 147 -> 185 days
$ws.Range("D43").Value = "The number of people who have a lung cancer diagnosis which is not associated with TLHC activity because they did not take up the offer of a LHC, or attended LHC but were assessed as low risk, or even had a scan but the diagnosis was made over 185 days following their scan."
$ws.Range("E43").Value = "The number of unique particiapnts in the invites table who have a lung cancer diagnosis in the NCRAS dataset but which is not associated with TLHC activity either because the participant didn't receive a scan (were invited but didn't take up the offer, attended LHC but were assessed as low risk) or did receive a scan but the diagnosis was made over 185 days afterwards."

# Row 44 (was 43) - 147 -> 185 days
$ws.Range("E44").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 185 days of a TLHC scan or LHC at which they were assessed as high risk and with a Tumour-Node-Metastasis (TNM) staging of either 1 or 2."

# Row 45 (was 44) - 147 -> 185 days
$ws.Range("E45").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 185 days of a TLHC scan or LHC at which they were assessed as high risk and with a Tumour-Node-Metastasis (TNM) staging of either 3 or 4."

# Row 46 (was 45) - 147 -> 185 days
$ws.Range("E46").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis was made within 185 days of a TLHC scan or LHC at which they were assessed as high risk and staging information is not provided because there is insufficient information or the cancer is unstageable."

# Row 47 (was 46) - 147 -> 185 days
$ws.Range("E47").Value = "The number of unique participants with a lung cancer diagnosis with a Tumour-Node-Metastasis (TNM) staging of either 1 or 2 and where the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 185 days following scan."

# Row 48 (was 47) - 147 -> 185 days
$ws.Range("E48").Value = "The number of unique participants with a lung cancer diagnosis with a Tumour-Node-Metastasis (TNM) staging of either 3 or 4 and where the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 185 days following scan."

# Row 49 (was 48) - 147 -> 185 days
$ws.Range("E49").Value = "The number of unique participants with a lung cancer diagnosis where the diagnosis is without staging information because there is insufficient information or the cancer is unstageable and the diagnosis is not associated with TLHC activity because the participant did not undergo a scan or the diagnosis was made over 185 days following scan."

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E49"))

Write-Host "Done"
